# Updated symbol list on Sat Feb 11 11:21:47 UTC 2023 with GitHub Actions
#
# This script refreshes the crypto-ranking snapshot on Sheet1:
#   - BNB/OKB/HuobiToken/Cronos price & 1h-volume figures are nudged to the
#     latest scrape
#   - GateToken jumped from rank 16 -> rank 4, so rows 6-17 (Coin name, Link,
#     Price, Volume(1h)) roll down by one slot, with GateToken's fresh row
#     landing at row 6
#   - Remaining rows keep their coin/link but get refreshed Price/Volume(1h)
#     values
#
# All Price/Volume(1h) cells are stored as plain text (e.g. "308.30",
# "0.28%") rather than numbers, so values are written with a leading
# apostrophe to stop Excel from auto-converting them to numeric/percentage
# types (which would silently drop things like trailing zeros).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: BNB ---
$ws.Range("D2").Value = "'308.30"
$ws.Range("E2").Value = "'0.28%"

# --- Row 3: OKB ---
$ws.Range("D3").Value = "'40.84"
$ws.Range("E3").Value = "'2.61%"

# --- Row 4: HuobiToken ---
$ws.Range("E4").Value = "'-0.32%"

# --- Row 5: Cronos ---
$ws.Range("D5").Value = "'0.07622"
$ws.Range("E5").Value = "'-1.47%"

# --- Row 6: FTXToken -> GateToken (new entrant pushes list down) ---
$ws.Range("B6").Value = "GateToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D6").Value = "'4.257"
$ws.Range("E6").Value = "'0.54%"

# --- Row 7: BTSEToken -> FTXToken ---
$ws.Range("B7").Value = "FTXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D7").Value = "'1.605"
$ws.Range("E7").Value = "'-0.60%"

# --- Row 8: MXToken -> BTSEToken ---
$ws.Range("B8").Value = "BTSEToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D8").Value = "'2.470"
$ws.Range("E8").Value = "'2.10%"

# --- Row 9: LiechtensteinCryptoassetsExchange -> MXToken ---
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D9").Value = "'0.9021"
$ws.Range("E9").Value = "'1.10%"

# --- Row 10: WazirX -> LiechtensteinCryptoassetsExchange ---
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").Value = "'0.1127"
$ws.Range("E10").Value = "'12.31%"

# --- Row 11: MandalaExchangeToken -> WazirX ---
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "'0.1790"
$ws.Range("E11").Value = "'3.05%"

# --- Row 12: BitrueCoin -> MandalaExchangeToken ---
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "'0.09113"
$ws.Range("E12").Value = "'1.19%"

# --- Row 13: BitMartToken -> BitrueCoin ---
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.04181"
$ws.Range("E13").Value = "'-5.85%"

# --- Row 14: BitForexToken -> BitMartToken ---
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.1053"
$ws.Range("E14").Value = "'-0.12%"

# --- Row 15: TigerCash -> BitForexToken ---
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "'0.001251"
$ws.Range("E15").Value = "'-0.43%"

# --- Row 16: LEO -> TigerCash ---
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "'0.005683"
$ws.Range("E16").Value = "'-2.78%"

# --- Row 17: GateToken -> LEO ---
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "'3.348"
$ws.Range("E17").Value = "'-0.14%"

# --- Row 18: BitpandaEcosystemToken ---
$ws.Range("E18").Value = "'-0.75%"

# --- Row 19: MCDex ---
$ws.Range("D19").Value = "'6.632"
$ws.Range("E19").Value = "'-5.97%"

# --- Row 20: ProBitToken ---
$ws.Range("E20").Value = "'1.30%"

# --- Row 21: ZBToken ---
$ws.Range("D21").Value = "'0.2801"
$ws.Range("E21").Value = "'1.54%"

# --- Row 22: CoinExToken ---
$ws.Range("D22").Value = "'0.04077"
$ws.Range("E22").Value = "'-1.51%"

# --- Row 23: BitKan ---
$ws.Range("D23").Value = "'0.001245"
$ws.Range("E23").Value = "'3.18%"

# --- Row 24: HotbitToken ---
$ws.Range("D24").Value = "'0.004093"
$ws.Range("E24").Value = "'0.72%"

# --- Row 25: NitroEx ---
$ws.Range("E25").Value = "'0.01%"

# --- Row 38: One ---
$ws.Range("D38").Value = "'0.02396"
$ws.Range("E38").Value = "'2.01%"

# --- Row 39: IDEX ---
$ws.Range("D39").Value = "'0.05176"
$ws.Range("E39").Value = "'-0.38%"

# --- Row 40: KickToken ---
$ws.Range("D40").Value = "'0.007777"
$ws.Range("E40").Value = "'-2.11%"

# --- Row 41: BKEXToken ---
$ws.Range("E41").Value = "'-1.71%"

# --- Row 42: Dexo ---
$ws.Range("D42").Value = "'0.007061"
$ws.Range("E42").Value = "'12.42%"

# --- Row 43: CEJI ---
$ws.Range("E43").Value = "'0.04%"

# --- Row 44: LocalTraders ---
$ws.Range("D44").Value = "'0.007709"
$ws.Range("E44").Value = "'-5.98%"

# --- Row 45: PooCoin ---
$ws.Range("D45").Value = "'0.3080"
$ws.Range("E45").Value = "'-7.34%"

# --- Row 46: CoinLion ---
$ws.Range("D46").Value = "'0.00006965"
$ws.Range("E46").Value = "'7.04%"

# --- Row 47: Kangarootoken ---
$ws.Range("E47").Value = "'0.01%"

# --- Row 48: BOLO ---
$ws.Range("D48").Value = "'0.04627"
$ws.Range("E48").Value = "'1,200.81%"

# --- Row 50: CryptobidCoin ---
$ws.Range("E50").Value = "'0.01%"

# --- Row 51: SpecialPowerGold ---
$ws.Range("E51").Value = "'0.01%"
